{"js": "// Update the court phone number referenced near the end of the letter:\n// \"can call the court on 01792 485 800.\" -> \"can call the court on 0330 808 4424.\"\nconst body = context.document.body;\n\nconst results = body.search(\"01792 485 800\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the phone number '01792 485 800' to replace.\");\n}\n\nfor (const rng of results.items) {\n  rng.insertText(\"0330 808 4424\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the court phone number referenced near the end of the letter:\n# \"You can call the court on 01792 485 800. The court cannot give legal advice.\"\n# -> \"You can call the court on 0330 808 4424. The court cannot give legal advice.\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n$find.Execute(\"01792 485 800\", $true, $false, $false, $false, $false, $true, 1, $false, \"0330 808 4424\", 2) | Out-Null\n"}
